{"js": "// Lattice-multiplication exercise table refresh.\n// The table is 5 rows x 3 columns = 15 cells; each cell holds a single\n// run with 5 text segments separated by manual line breaks:\n//   \"A x B\" / \"  d1    d2\" (digits of B spaced out) / \"  ----\" /\n//   \"<digit0 of A>|    |\" / \"<digit1 of A>|    |\"\n// This script overwrites the \"A x B\" pair for every cell (in row-major,\n// left-to-right/top-to-bottom order) and regenerates the other four\n// lines from the new A/B values, leaving all formatting untouched.\n\nconst newPairs = [\n  [\"22\", \"86\"],\n  [\"39\", \"89\"],\n  [\"13\", \"26\"],\n  [\"73\", \"54\"],\n  [\"78\", \"68\"],\n  [\"31\", \"41\"],\n  [\"14\", \"80\"],\n  [\"83\", \"61\"],\n  [\"81\", \"96\"],\n  [\"53\", \"83\"],\n  [\"81\", \"75\"],\n  [\"66\", \"47\"],\n  [\"82\", \"71\"],\n  [\"56\", \"84\"],\n  [\"24\", \"47\"],\n];\n\nconst LINE_BREAK = \"\\u000b\"; // Word manual line break (<w:br/>) in Office.js text\n\nfunction buildCellText(a, b) {\n  const line1 = a + \" x \" + b;\n  const line2 = \"  \" + b.split(\"\").join(\"    \");\n  const line3 = \"  ----\";\n  const line4 = a[0] + \"|    |\";\n  const line5 = a[1] + \"|    |\";\n  return [line1, line2, line3, line4, line5].join(LINE_BREAK);\n}\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet pairIndex = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const cellBody = cell.body;\n    cellBody.paragraphs.load(\"items\");\n    await context.sync();\n\n    const [a, b] = newPairs[pairIndex];\n    pairIndex++;\n\n    const paragraph = cellBody.paragraphs.items[0];\n    paragraph.insertText(buildCellText(a, b), \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication exercise table refresh.\n# The table is 5 rows x 3 columns = 15 cells; each cell holds a single\n# run with 5 text segments separated by manual line breaks:\n#   \"A x B\" / \"  d1    d2\" (digits of B spaced out) / \"  ----\" /\n#   \"<digit0 of A>|    |\" / \"<digit1 of A>|    |\"\n# This script overwrites the \"A x B\" pair for every cell (in row-major,\n# left-to-right/top-to-bottom order) and regenerates the other four\n# lines from the new A/B values, leaving all formatting untouched.\n\n$newPairs = @(\n    @(\"22\", \"86\"),\n    @(\"39\", \"89\"),\n    @(\"13\", \"26\"),\n    @(\"73\", \"54\"),\n    @(\"78\", \"68\"),\n    @(\"31\", \"41\"),\n    @(\"14\", \"80\"),\n    @(\"83\", \"61\"),\n    @(\"81\", \"96\"),\n    @(\"53\", \"83\"),\n    @(\"81\", \"75\"),\n    @(\"66\", \"47\"),\n    @(\"82\", \"71\"),\n    @(\"56\", \"84\"),\n    @(\"24\", \"47\")\n)\n\n$LINE_BREAK = [char]11  # Word manual line break (<w:br/>)\n\nfunction Build-CellText($a, $b) {\n    $line1 = \"$a x $b\"\n    $line2 = \"  \" + ($b.ToCharArray() -join \"    \")\n    $line3 = \"  ----\"\n    $line4 = \"$($a.Substring(0,1))|    |\"\n    $line5 = \"$($a.Substring(1,1))|    |\"\n    return $line1 + $LINE_BREAK + $line2 + $LINE_BREAK + $line3 + $LINE_BREAK + $line4 + $LINE_BREAK + $line5\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$index = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $newPairs[$index]\n        $index++\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = Build-CellText $pair[0] $pair[1]\n    }\n}\n"}
